$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.319.44"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "'1.933.93"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -1.21%  "
$ws.Range("D5").Value = "'251.38"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "'0.7153"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").Value = "'0.3297"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "'27.64"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("D10").Value = "'0.07253"
$ws.Range("E10").Value = "  +5.58%  "
$ws.Range("D11").Value = "'0.8022"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "'0.08100"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "'1.930.38"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'5.473"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "'94.75"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "'15.01"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").Value = "'30.310.14"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "'252.62"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("D19").Value = "'0.000008184"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("D20").Value = "'5.814"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "'2.185.37"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "'6.959"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "'9.747"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'164.79"
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("D27").Value = "'2.352"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").Value = "'19.33"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "'0.1301"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "'1.353"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").Value = "'1.538"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "'4.428"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "'4.177"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").Value = "'0.05211"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "'1.263"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").Value = "'0.7469"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").Value = "'2.777"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "'0.01968"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "'78.83"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("D41").Value = "'6.423"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("D42").Value = "'0.4531"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "'2.024"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").Value = "'0.8424"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").Value = "'101.47"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").Value = "'9.767"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").Value = "'7.431"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "'36.79"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "'0.4176"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'0.06034"
$ws.Range("E51").Value = "  +0.82%  "
